$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H holds the "Save" metric, add header to match the style
# of the existing header row (bold, bordered, centered - same as G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "Save"

# Populate the new column's data rows with the Save values.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
